$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1976
$ws.Range("J3").Value = 8080
$ws.Range("K3").Value = 1902
$ws.Range("J4").Value = 1803
$ws.Range("K4").Value = 403
$ws.Range("K5").Value = 125
$ws.Range("K6").Value = 2428
$ws.Range("J7").Value = 29274
$ws.Range("K7").Value = 6834

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 131
$ws.Range("J3").Value = 530
$ws.Range("K3").Value = 132
$ws.Range("K6").Value = 161
$ws.Range("J7").Value = 1852
$ws.Range("K7").Value = 459

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J4").Value = 23
$ws.Range("K6").Value = 31
$ws.Range("J7").Value = 591
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 76
$ws.Range("K3").Value = 107
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 277

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 41
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 40
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 200
$ws.Range("J8").Value = 1852
$ws.Range("K8").Value = 459
$ws.Range("K10").Value = 40
$ws.Range("K14").Value = 37
$ws.Range("K18").Value = 48
$ws.Range("K19").Value = 188
$ws.Range("K20").Value = 148
$ws.Range("K23").Value = 62
$ws.Range("K25").Value = 32
$ws.Range("K27").Value = 78
$ws.Range("K29").Value = 335
$ws.Range("K33").Value = 277
$ws.Range("K34").Value = 41
$ws.Range("K37").Value = 229
$ws.Range("K41").Value = 62
$ws.Range("K42").Value = 238
$ws.Range("K48").Value = 81
$ws.Range("K49").Value = 40
$ws.Range("K50").Value = 38
$ws.Range("K51").Value = 78
$ws.Range("K52").Value = 185
$ws.Range("K54").Value = 121
$ws.Range("K60").Value = 47
$ws.Range("J63").Value = 96
$ws.Range("K63").Value = 23
$ws.Range("K65").Value = 165
$ws.Range("K67").Value = 264
$ws.Range("K72").Value = 30
$ws.Range("K77").Value = 47
$ws.Range("K79").Value = 180
$ws.Range("J83").Value = 591
$ws.Range("K83").Value = 147
$ws.Range("K85").Value = 344
$ws.Range("K88").Value = 90
$ws.Range("K89").Value = 90
$ws.Range("K91").Value = 65
$ws.Range("K95").Value = 109
$ws.Range("K98").Value = 46
$ws.Range("K99").Value = 123
$ws.Range("J101").Value = 29274
$ws.Range("K101").Value = 6834

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 75
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 42
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 110
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 67
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 20
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 64
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K5").Value = 9
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 200

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 20
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 116
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 344

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 40
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 185
